# Updates the "ESTADO DE CUENTA" worker list (rows 16-24) on Hoja1.
#
# The underlying data set changed: a new pending period (2108) was added for
# MARIA ANGELICA PEREZ HERNANDEZ, a brand-new worker (JULIETH PAOLA
# VILLARREAL MARTINEZ) was inserted, ELOINA MARIA SARABIA SIMARRA's
# "Salario Basico" value was corrected, and RAISA CONEO CAMERO's three
# overdue periods were re-ordered chronologically (2308, 2309, 2310).
#
# Net effect: rows 16-24, columns C (N Doc Trabajador), D (Nombre
# Trabajador), E (Periodo Mora), F (Valor Mora) and G (Salario Basico) get
# rewritten with the refreshed data below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 16; C = "1049453756"; D = "MARIA ANGELICA PEREZ HERNANDEZ";     E = "2108"; F = 36341; G = 908526  },
    @{ Row = 17; C = "1049453756"; D = "MARIA ANGELICA PEREZ HERNANDEZ";     E = "2109"; F = 36341; G = 908526  },
    @{ Row = 18; C = "45691997";   D = "RUTH MARIA SCHLEGEL CORREA";         E = "2111"; F = 36341; G = 908526  },
    @{ Row = 19; C = "22807446";   D = "KATTY HEIDY GUERRERO CAVANA";        E = "2111"; F = 60000; G = 1500000 },
    @{ Row = 20; C = "1047457178"; D = "ELOINA MARIA SARABIA SIMARRA";       E = "2111"; F = 36341; G = 908526  },
    @{ Row = 21; C = "1047499055"; D = "JULIETH PAOLA VILLARREAL MARTINEZ";  E = "2111"; F = 36341; G = 908526  },
    @{ Row = 22; C = "1047426149"; D = "RAISA CONEO CAMERO";                 E = "2308"; F = 60000; G = 1500000 },
    @{ Row = 23; C = "1047426149"; D = "RAISA CONEO CAMERO";                 E = "2309"; F = 60000; G = 1500000 },
    @{ Row = 24; C = "1047426149"; D = "RAISA CONEO CAMERO";                 E = "2310"; F = 60000; G = 1500000 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
}
